# Update cryptocurrency symbol list (prices + volume labels) to the
# latest scrape snapshot, as produced by the scheduled GitHub Actions run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores numeric-looking values as TEXT (the sheet
# is generated by a scraper, not typed by hand), e.g. "22.80" must keep
# its trailing zero and must not turn into the number 22.8. Assigning a
# numeric-looking string straight to .Value on a General-formatted cell
# makes Excel auto-convert it to a real number, so we momentarily force
# the cell to Text format, assign the literal string, then restore the
# cell's original style so no stray formatting is left behind.
function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = $origStyle
}

Set-TextValue "D2"  "248.01"
Set-TextValue "D3"  "22.80"
Set-TextValue "D4"  "5.301"
Set-TextValue "D5"  "0.05736"
Set-TextValue "D6"  "3.433"
Set-TextValue "D7"  "0.8091"
Set-TextValue "D8"  "0.8717"
Set-TextValue "D10" "0.07376"
Set-TextValue "D12" "0.03112"
Set-TextValue "D13" "0.09397"
Set-TextValue "D14" "3.893"
Set-TextValue "D15" "0.001578"
Set-TextValue "D16" "0.04816"

# Column E (Volume(1h)) holds plain text labels, e.g. "16OneONEWorstin24h"
# becoming "16OneONE" -- no numeric coercion risk, so assign directly.
$ws.Range("E17").Value = "16OneONE"

Set-TextValue "D18" "0.006142"
Set-TextValue "D19" "0.005164"
Set-TextValue "D20" "0.0009976"
Set-TextValue "D21" "0.0001500"
Set-TextValue "D22" "3.724"
Set-TextValue "D23" "6.324"
Set-TextValue "D24" "2.183"
Set-TextValue "D25" "0.3280"

Set-TextValue "D40" "0.03945"
Set-TextValue "D41" "0.006746"
Set-TextValue "D42" "0.1070"
Set-TextValue "D43" "0.002210"
Set-TextValue "D44" "0.007278"
Set-TextValue "D45" "0.00005606"

Set-TextValue "D47" "0.6001"

Set-TextValue "D48" "0.1778"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"

Set-TextValue "D49" "0.00002100"
